$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Link"
$ws.Range("B2").Value = "https://hi-viewer.web.app/mirador/?manifest=https://hi-ut.github.io/dataset/iiif/collection/nishikie_hi.json"
$ws.Range("B3").Value = "https://hi-viewer.web.app/mirador/?manifest=https://hi-ut.github.io/dataset/iiif/collection/nishikie_shizuoka.json"
$ws.Range("B4").Value = "https://hi-viewer.web.app/mirador/?manifest=https://hi-ut.github.io/dataset/iiif/collection/nishikie_yokohama.json"
